# Adds 2 new test cases (3 data rows) to the "Test Cases" sheet, mirroring
# the formatting of existing rows as closely as COM allows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# A cell elsewhere in the workbook that already carries the "s=1" style
# (plain font + thin border, no wrap) that the newly-authored rows use.
$plainBorderSrc = $wb.Worksheets.Item("AuthoringTest").Range("A2")

$rows = @(
    @{
        Row = 42
        TCID = "VerifyPostRecordDetails"
        Jira = "OPQA-370"
        Desc = "Verify that user contributed articles display the information about the author"
    },
    @{
        Row = 43
        TCID = "SeacrhAndViewOwnPost"
        Jira = "OPQA-415"
        Desc = "Verify that user is able to search the  posts a user authored themselves and view them."
    },
    @{
        Row = 44
        TCID = "SeacrhAndViewOthersPost"
        Jira = "OPQA-416"
        Desc = "Verify that user is able to search the posts of others and view them."
    }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $cellA = $ws.Cells.Item($rowNum, 1)
    $cellB = $ws.Cells.Item($rowNum, 2)
    $cellC = $ws.Cells.Item($rowNum, 3)
    $cellD = $ws.Cells.Item($rowNum, 4)
    $cellE = $ws.Cells.Item($rowNum, 5)

    # Match formatting of the previous data row for column B (hyperlink-ish font+border),
    # and reuse the plain bordered style (no wrap) used by these newly authored rows
    # for the other columns.
    $ws.Cells.Item($rowNum - 1, 2).Copy()
    $cellB.PasteSpecial(-4122)

    $plainBorderSrc.Copy()
    $cellA.PasteSpecial(-4122)
    $cellC.PasteSpecial(-4122)
    $cellD.PasteSpecial(-4122)
    $cellE.PasteSpecial(-4122)

    $cellA.Value = $r.TCID
    $cellB.Value = $r.Jira
    $cellC.Value = $r.Desc
    $cellD.Value = "Y"
    $cellE.Value = "PASS"
}

$ws.Range("B44").Select()

# Scroll the sheet so row 28 is the first visible row (best-effort: mirrors
# the authored workbook's `topLeftCell="A28"`).
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 90
